$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.150.08'
$ws.Range('E2').Value = '  +3.09%  '
$ws.Range('D3').Value = '2.066.74'
$ws.Range('E3').Value = '  +2.69%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '230.60'
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.616'
$ws.Range('E6').Value = '  +1.49%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '58.17'
$ws.Range('E7').Value = '  +6.18%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  +2.08%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0808'
$ws.Range('E10').Value = '  +2.40%  '
$ws.Range('E11').Value = '  -0.84%  '
$ws.Range('D12').Value = '2.381.38'
$ws.Range('E12').Value = '  +3.07%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.66'
$ws.Range('E13').Value = '  +2.93%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.75'
$ws.Range('E14').Value = '  +2.06%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.751'
$ws.Range('E15').Value = '  +1.41%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.28'
$ws.Range('E16').Value = '  +2.97%  '
$ws.Range('D17').Value = '2.078.22'
$ws.Range('E17').Value = '  +3.57%  '
$ws.Range('D18').Value = '38.123.28'
$ws.Range('E18').Value = '  +3.27%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.22'
$ws.Range('E19').Value = '  +1.89%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.93'
$ws.Range('E20').Value = '  +1.70%  '
$ws.Range('D21').Value = '0.0₃0833'
$ws.Range('E21').Value = '  +1.46%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '224.87'
$ws.Range('E22').Value = '  +0.44%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  +0.89%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.26'
$ws.Range('E25').Value = '  +2.37%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.27'
$ws.Range('E26').Value = '  +0.90%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.66'
$ws.Range('E27').Value = '  +0.38%  '
$ws.Range('E28').Value = '  +6.90%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.08'
$ws.Range('E29').Value = '  +2.03%  '
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('E31').Value = '  +1.66%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.56'
$ws.Range('E32').Value = '  +0.39%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.61'
$ws.Range('E33').Value = '  +3.89%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').Value = '  +7.64%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.39'
$ws.Range('E36').Value = '  +1.98%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.01'
$ws.Range('E37').Value = '  +11.49%  '
$ws.Range('E38').Value = '  +5.13%  '
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '98.49'
$ws.Range('E40').Value = '  +3.49%  '
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('D42').Value = '1.483.54'
$ws.Range('E42').Value = '  +0.44%  '
$ws.Range('E43').Value = '  +3.04%  '
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.89'
$ws.Range('E44').Value = '  +4.66%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '16.86'
$ws.Range('E45').Value = '  +2.11%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.19'
$ws.Range('E46').Value = '  +20.08%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.13'
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('E48').Value = '  +2.06%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.96'
$ws.Range('E49').Value = '  +2.08%  '
$ws.Range('E50').Value = '  -1.65%  '
$ws.Range('D51').Value = '2.265.96'
$ws.Range('E51').Value = '  +2.85%  '
